# Update functional requirement v1.1
# Rewrites the "Functional Requirements" sheet (rows 5-14) with the
# reorganised/updated requirement list, keeping FR-11/FR-12/FR-13 rows
# (15-17) untouched content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional Requirements")

# --- Row 5: FR-01 / Inventory management ---
$ws.Range("B5").Value = "FR-01"
$ws.Range("C5").Value = "Inventory management"
$ws.Range("D5").Value = "Must"
$ws.Range("E5").Value = "Accountants shall be able to  manage inventory by creating a Goods Received Note when goods are imported into the warehouse. "
$ws.Rows.Item(5).RowHeight = 45

# --- Row 6: FR-02 / Order management ---
$ws.Range("B6").Value = "FR-02"
$ws.Range("C6").Value = "Order management"
$ws.Range("D6").Value = "Must"
$ws.Range("E6").Value = "Accountants shall be able to process orders, track and update orders and payments status, generate invoices.  Resellers/customers should also be able to view the payment status of their orders."
$ws.Rows.Item(6).RowHeight = 60

# --- Row 7: FR-03 / Order Placement ---
$ws.Range("B7").Value = "FR-03"
$ws.Range("C7").Value = "Order Placement"
$ws.Range("D7").Value = "Must"
$ws.Range("E7").Value = "Resellers/customers should be able to place an order for items by selecting the desired products, specifying the quantity, and choosing a payment method (Cash, bank transfer, Momo...)."
$ws.Rows.Item(7).RowHeight = 60

# --- Row 8: FR-04 / Order Status Tracking ---
$ws.Range("B8").Value = "FR-04"
$ws.Range("C8").Value = "Order Status Tracking"
$ws.Range("D8").Value = "Must"
$ws.Range("E8").Value = "Resellers/customers should be able to track the status of their orders, including knowing when the order has been processed, shipped, or delivered."
$ws.Rows.Item(8).RowHeight = 45

# --- Row 9: FR-05 / User Management ---
$ws.Range("B9").Value = "FR-05"
$ws.Range("C9").Value = "User Management"
$ws.Range("D9").Value = "Must"
$ws.Range("E9").Value = "Operators shall be able to manage user profiles and accounts, assign user roles and permissions, control access to sensitive data, generating reports on users, …"
$ws.Rows.Item(9).RowHeight = 60

# --- Row 10: FR-06 / Product Management ---
$ws.Range("B10").Value = "FR-06"
$ws.Range("C10").Value = "Product Management"
$ws.Range("D10").Value = "Must"
$ws.Range("E10").Value = "Accountants shall be able to add, delete, edit product information."
$ws.Rows.Item(10).RowHeight = 30

# --- Row 11: FR-07 / Delivery Management ---
# This row previously carried the special fill/border highlight; in the
# updated sheet it reverts to the plain body style, so grab the (already
# plain) format from row 10 before writing the new text.
$ws.Range("C10:D10").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)
$ws.Range("B11").Value = "FR-07"
$ws.Range("C11").Value = "Delivery Management"
$ws.Range("D11").Value = "Must"
$ws.Range("E11").Value = "Accountants shall be able to create goods delivery note to deliver goods to resellers, update the status of orders as being transferred."
$ws.Rows.Item(11).RowHeight = 45

# --- Row 12: FR-08 / Payment Integration (keeps its highlighted style) ---
$ws.Range("B12").Value = "FR-08"
$ws.Range("C12").Value = "Payment Integration"
$ws.Range("D12").Value = "Must"
$ws.Range("E12").Value = "The software should be integrated with different payment gateways, allowing resellers/customers to make online payments through the platform."
$ws.Rows.Item(12).RowHeight = 45

# --- Row 13: FR-09 / Stock Reporting (keeps its highlighted style) ---
$ws.Range("B13").Value = "FR-09"
$ws.Range("C13").Value = "Stock Reporting"
$ws.Range("D13").Value = "Must"
$ws.Range("E13").Value = "The software should provide incoming/outgoing stock reports, allowing accountants to view stock levels, track inventory movements, and manage stock levels efficiently."
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: FR-10 / Sales Reporting ---
$ws.Range("B14").Value = "FR-10"
$ws.Range("C14").Value = "Sales Reporting"
$ws.Range("D14").Value = "Must"
$ws.Range("E14").Value = "The software should provide sales reporting and revenue reports, allowing accountants to view best-selling products, track sales trends, and forecast sales performance."
$ws.Rows.Item(14).RowHeight = 60

# Update the view: selection moved to D20, no pinned top-left scroll cell.
$ws.Range("D20").Select()
